$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Answers")
$ws.Range("B1").Value = "Key Answers"
$ws.Range("B1").Select()
